$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.202.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").Value = "'2.372.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'501.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").Value = "'130.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.04%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").Value = "'2.375.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.18%  "
$ws.Range("D10").Value = "'0.0982"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "'4.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "'2.796.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "'56.176.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").Value = "'21.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "'2.407.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.63%  "
$ws.Range("D19").Value = "'10.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.98%  "
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("D21").Value = "'306.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").Value = "'6.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.90%  "
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("D24").Value = "'64.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "'0.369"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.93%  "
$ws.Range("E27").Value = "  -4.59%  "
$ws.Range("D28").Value = "'7.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.36%  "
$ws.Range("D29").Value = "'172.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("E31").Value = "  -3.69%  "
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").Value = "'5.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.58%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.00%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "'0.996"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").Value = "'17.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "  -6.05%  "
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("D39").Value = "'35.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("E41").Value = "  -4.24%  "
$ws.Range("D42").Value = "'130.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.14%  "
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("E44").Value = "  -5.17%  "
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").Value = "'0.0902"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").Value = "'241.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.37%  "
$ws.Range("E48").Value = "  -2.98%  "
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("E51").Value = "  -3.45%  "
